# "New Unlocks and Requirements for Page and more"
#
# 1. Rename Sheet5 -> Enemies
# 2. Fix the misspelled creature name on the (now) Enemies sheet
# 3. Add the new Mono/Dual unlock-requirement legend to Sheet1 (rows 25-29)
# 4. Restore the recorded selections / active sheet for each tab

$wb = $excel.ActiveWorkbook

# --- 1 & 2: Enemies sheet -------------------------------------------------
$wsEnemies = $wb.Worksheets.Item("Sheet5")
$wsEnemies.Name = "Enemies"
$wsEnemies.Range("B2").Value = "bringobrongo"

# --- 3: New legend rows on Sheet1 ------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("D25").Value = "Mono:"
$ws1.Range("E25").Value = "Unlock"
$ws1.Range("F25").Value = "Upgrade"

$ws1.Range("D26").Value = "All"
$ws1.Range("E26").Value = "Default"
$ws1.Range("F26").Value = "Level"

$ws1.Range("D28").Value = "Dual:"

$ws1.Range("D29").Value = "Brute"
$ws1.Range("E29").Value = "15 CW"

# --- 4: Recorded selections -------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("D27").Select()

$wsEnemies.Range("A3").Select()

# Sheet1 is the active tab when the workbook was saved, so activate it last.
$ws1.Activate()
$ws1.Range("K22").Select()
